$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, pushing existing rows 75-82 down to 76-83.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly price record.
$ws.Range("A75").Value = 11
$ws.Range("B75").Value = "Vega Monumental Concepción"
$ws.Range("C75").Value = "Bíobío"
$ws.Range("D75").Value = 45275
$ws.Range("E75").Value = 8
$ws.Range("F75").Value = 100112026
$ws.Range("G75").Value = "Haba"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 80
$ws.Range("K75").Value = 19000
$ws.Range("L75").Value = 19000
$ws.Range("M75").Value = 19000
$ws.Range("N75").Value = "$/saco 25 kilos"
$ws.Range("O75").Value = "Región del Maule"
$ws.Range("P75").Value = 760
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = "Hortaliza"
